# Apply the 2022-12-01 YTD data refresh across Citywide Totals, By Neighborhood,
# and the affected individual neighborhood sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("F3").Value = 135
$ws.Range("F4").Value = 8
$ws.Range("H6").Value = 436
$ws.Range("C6").Value = 472
$ws.Range("F6").Value = 520
$ws.Range("E6").Value = 463
$ws.Range("B6").Value = 370
$ws.Range("D6").Value = 410
$ws.Range("G6").Value = 433
$ws.Range("I6").Value = 496
$ws.Range("E7").Value = 686
$ws.Range("G7").Value = 662
$ws.Range("B7").Value = 495
$ws.Range("I7").Value = 826
$ws.Range("F7").Value = 755
$ws.Range("C7").Value = 626
$ws.Range("H7").Value = 706
$ws.Range("D7").Value = 639

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("E6").Value = 52
$ws.Range("E7").Value = 65

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("D6").Value = 21
$ws.Range("D7").Value = 36

$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range("F3").Value = 1
$ws.Range("F6").Value = 14

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("I6").Value = 21
$ws.Range("H6").Value = 26
$ws.Range("H7").Value = 45
$ws.Range("I7").Value = 44

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("F5").Value = 14
$ws.Range("B8").Value = 30
$ws.Range("F8").Value = 50
$ws.Range("B16").Value = 2
$ws.Range("D19").Value = 27
$ws.Range("H28").Value = 45
$ws.Range("I28").Value = 44
$ws.Range("F29").Value = 13
$ws.Range("E32").Value = 65
$ws.Range("D36").Value = 36
$ws.Range("F47").Value = 17
$ws.Range("F51").Value = 7
$ws.Range("F53").Value = 80
$ws.Range("I53").Value = 123
$ws.Range("D61").Value = 3
$ws.Range("E62").Value = 7
$ws.Range("C65").Value = 22
$ws.Range("F70").Value = 24
$ws.Range("I75").Value = 2
$ws.Range("F76").Value = 19
$ws.Range("G77").Value = 24
$ws.Range("G97").Value = 5
$ws.Range("E98").Value = 686
$ws.Range("G98").Value = 662
$ws.Range("H98").Value = 706
$ws.Range("D98").Value = 639
$ws.Range("I98").Value = 826
$ws.Range("B98").Value = 495
$ws.Range("C98").Value = 626
$ws.Range("F98").Value = 755

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("F6").Value = 59
$ws.Range("I6").Value = 78
$ws.Range("F7").Value = 80
$ws.Range("I7").Value = 123

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("F4").Value = 2
$ws.Range("F6").Value = 19

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("C5").Value = 18
$ws.Range("C6").Value = 22

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("E4").Value = 5
$ws.Range("E5").Value = 7

$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range("F5").Value = 12
$ws.Range("F6").Value = 13

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("G6").Value = 14
$ws.Range("G7").Value = 24

$ws = $wb.Worksheets.Item('Wrigleyville')
$ws.Range("G5").Value = 3
$ws.Range("G6").Value = 5

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("D5").Value = 15
$ws.Range("D6").Value = 27

$ws = $wb.Worksheets.Item('Bucktown')
$ws.Range("B4").Value = 2
$ws.Range("B5").Value = 2

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("F4").Value = 6
$ws.Range("F5").Value = 7

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("F5").Value = 10
$ws.Range("F6").Value = 17

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("F4").Value = 18
$ws.Range("F5").Value = 24

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range("H4").Value = 1
$ws.Range("H5").Value = 2

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("F3").Value = 6
$ws.Range("B5").Value = 21
$ws.Range("B6").Value = 30
$ws.Range("F6").Value = 50
